$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of H2 (removing the extra shared-string reference)
$ws.Range("H2").ClearContents()

# Move/update the current selection to H2
$ws.Range("H2").Select()
